$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/3/2025  Through  2/9/2025"

# --- Weekly crime statistics table updates (rows 15-30) ---

# Cells that must hold the literal text "0" (not the number 0) -- format as Text first
# so the digit string is not auto-coerced into a numeric value.
$textZeroCells = @("D15","D23","D27","C29","D29","C30","D30")
foreach ($ref in $textZeroCells) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = "0"
}

# Cells that hold the literal text "***.*" (non-numeric, assigns as text naturally)
$ws.Range("E15").Value = "***.*"
$ws.Range("E23").Value = "***.*"
$ws.Range("E27").Value = "***.*"
$ws.Range("E29").Value = "***.*"
$ws.Range("E30").Value = "***.*"

# Numeric cell updates
$ws.Range("F15").Value = 6
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 9
$ws.Range("K15").Value = 125
$ws.Range("L15").Value = 125
$ws.Range("M15").Value = 12.5
$ws.Range("N15").Value = 50
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 14.285714285714
$ws.Range("F16").Value = 36
$ws.Range("G16").Value = 37
$ws.Range("H16").Value = -2.702702702702
$ws.Range("I16").Value = 49
$ws.Range("J16").Value = 54
$ws.Range("K16").Value = -9.259259259259
$ws.Range("L16").Value = 8.888888888888
$ws.Range("M16").Value = -12.5
$ws.Range("N16").Value = -68.181818181818
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 57.142857142857
$ws.Range("F17").Value = 52
$ws.Range("G17").Value = 60
$ws.Range("H17").Value = -13.333333333333
$ws.Range("I17").Value = 78
$ws.Range("J17").Value = 87
$ws.Range("K17").Value = -10.344827586206
$ws.Range("L17").Value = 18.181818181818
$ws.Range("M17").Value = 105.263157894737
$ws.Range("N17").Value = 4
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -22.727272727272
$ws.Range("I18").Value = 22
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = -24.137931034482
$ws.Range("L18").Value = -33.333333333333
$ws.Range("M18").Value = -40.54054054054
$ws.Range("N18").Value = -90.308370044052
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = 11.764705882352
$ws.Range("F19").Value = 74
$ws.Range("G19").Value = 86
$ws.Range("H19").Value = -13.953488372093
$ws.Range("I19").Value = 102
$ws.Range("J19").Value = 112
$ws.Range("K19").Value = -8.928571428571
$ws.Range("L19").Value = 67.213114754098
$ws.Range("M19").Value = 325
$ws.Range("N19").Value = 100
$ws.Range("C20").Value = 21
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 133.333333333333
$ws.Range("F20").Value = 63
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 79
$ws.Range("J20").Value = 56
$ws.Range("K20").Value = 41.071428571428
$ws.Range("L20").Value = -7.058823529411
$ws.Range("M20").Value = 51.923076923076
$ws.Range("N20").Value = -54.06976744186
$ws.Range("C21").Value = 64
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = 42.222222222222
$ws.Range("F21").Value = 248
$ws.Range("G21").Value = 251
$ws.Range("H21").Value = -1.195219123505
$ws.Range("I21").Value = 339
$ws.Range("J21").Value = 343
$ws.Range("K21").Value = -1.166180758017
$ws.Range("L21").Value = 14.915254237288
$ws.Range("M21").Value = 57.674418604651
$ws.Range("N21").Value = -50.583090379008
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 33.333333333333
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 33.333333333333
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 37.5
$ws.Range("I23").Value = 14
$ws.Range("K23").Value = 16.666666666666
$ws.Range("L23").Value = 7.692307692307
$ws.Range("M23").Value = 366.666666666667
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 4.545454545454
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = -0.934579439252
$ws.Range("I24").Value = 146
$ws.Range("J24").Value = 141
$ws.Range("K24").Value = 3.54609929078
$ws.Range("L24").Value = 24.786324786324
$ws.Range("M24").Value = 124.615384615385
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = -20.588235294117
$ws.Range("I25").Value = 34
$ws.Range("J25").Value = 40
$ws.Range("K25").Value = -15
$ws.Range("L25").Value = -2.857142857142
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 75
$ws.Range("G26").Value = 79
$ws.Range("H26").Value = -5.06329113924
$ws.Range("I26").Value = 108
$ws.Range("J26").Value = 103
$ws.Range("K26").Value = 4.854368932038
$ws.Range("L26").Value = 42.105263157894
$ws.Range("M26").Value = 27.058823529411
$ws.Range("F27").Value = 7
$ws.Range("H27").Value = 40
$ws.Range("I27").Value = 10
$ws.Range("K27").Value = 66.666666666666
$ws.Range("L27").Value = 42.857142857142
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 75
$ws.Range("I28").Value = 9
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = 28.571428571428
$ws.Range("L28").Value = -18.181818181818
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = -66.666666666666
$ws.Range("N29").Value = -50
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -33.333333333333
$ws.Range("N30").Value = -57.142857142857
